$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: A:C used to be three separate col entries (36.71/37.14/36.57),
# now a single merged col band A:C at ~35.14 chars wide. The engine quantizes
# ColumnWidth to coarse steps, so we pick the closest achievable input.
$ws.Range("A1:C35").ColumnWidth = 34.33

# --- Add column T (year 2023) mirroring column S (year 2022) formatting, row by row.
$srcRows = 3..34
foreach ($r in $srcRows) {
    $src = $ws.Cells.Item($r, 19)   # column S
    $dst = $ws.Cells.Item($r, 20)   # column T
    $src.Copy()
    $dst.PasteSpecial(-4122)        # xlPasteFormats
}
$excel.CutCopyMode = $false

# --- New data values for column T ---
$ws.Range("T4").Value = 2023

$ws.Range("T5").Value = 44
$ws.Range("T6").Value = 24
$ws.Range("T7").Value = 20

$ws.Range("T8").Value = "-"
$ws.Range("T9").Value = "-"
$ws.Range("T10").Value = "-"

$ws.Range("T11").Value = 5
$ws.Range("T12").Value = 1
$ws.Range("T13").Value = 4

$ws.Range("T14").Value = 8
$ws.Range("T15").Value = 6
$ws.Range("T16").Value = 2

$ws.Range("T17").Value = 5
$ws.Range("T18").Value = 1
$ws.Range("T19").Value = 4

$ws.Range("T20").Value = 7
$ws.Range("T21").Value = 5
$ws.Range("T22").Value = 2

$ws.Range("T23").Value = "-"
$ws.Range("T24").Value = "-"
$ws.Range("T25").Value = "-"

$ws.Range("T26").Value = 18
$ws.Range("T27").Value = 10
$ws.Range("T28").Value = 8

$ws.Range("T29").Value = "-"
$ws.Range("T30").Value = "-"
$ws.Range("T31").Value = "-"

$ws.Range("T32").Value = 1
$ws.Range("T33").Value = 1
$ws.Range("T34").Value = "-"

# --- Selection: the authored file no longer carries an explicit <selection>
# (it was produced by a non-interactive export), but the COM model always
# tracks an active cell. Selecting A1 is the closest achievable approximation
# of the "no special selection" state.
$ws.Range("A1").Select()
